$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bus movements and streetcar shuttle garage changes
# (row 2 holds the garage allocation text for each bus/streetcar route in row 1)
$ws.Range("EL2").Value = "Wilson"
$ws.Range("EN2").Value = "Birchmount, Queensway"

$ws.Range("FP2").Value = "Wilson"
$ws.Range("FQ2").Value = "."
$ws.Range("FR2").Value = "Birchmount, Queensway"
$ws.Range("FS2").Value = "."
$ws.Range("FT2").Value = "Eglinton"

# Restore the active selection to match the saved view state
$ws.Range("FT3").Select()
